# 2024.01.16 每日一题 hard cv
# Add two new LeetCode tracking rows (734. 句子相似性 / 2719. 统计整数数目)
# to the bottom of the table on Sheet1, with matching hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: 734. 句子相似性 --------------------------------------------
# Clone row 30's formatting (style + row height) for columns A:E, then
# overwrite the values so the new row matches the look of the existing table.
$ws.Range("A30:E30").Copy($ws.Range("A31:E31"))

$ws.Range("A31").Value = "734. 句子相似性"
$ws.Range("B31").Value = "哈希表"
$ws.Range("C31").Value = "⭐⭐"
$ws.Range("D31").Value = "2024.01.16"
$ws.Range("E31").Value = "力扣官方"

# TextToDisplay mirrors the cell text into the hyperlink's `display=`
# attribute (as the A-column links elsewhere in the sheet have), but it also
# overwrites the cell's own text with that string - so restore the intended
# title afterwards. The E-column "参考题解" links never carry a display
# attribute in this sheet, so leave TextToDisplay unset for those.
$ws.Hyperlinks.Add($ws.Range("A31"), "https://leetcode.cn/problems/sentence-similarity/", "", "", "https://leetcode.cn/problems/sentence-similarity/") | Out-Null
$ws.Range("A31").Value = "734. 句子相似性"
$ws.Hyperlinks.Add($ws.Range("E31"), "https://leetcode.cn/problems/sentence-similarity/solutions/2612921/ju-zi-xiang-si-xing-by-leetcode-solution/?envType=daily-question&envId=2024-01-16") | Out-Null

# Adding a hyperlink nudges Excel into minting a fresh "hyperlink" cell
# style; re-stamp the original column formatting so A31/E31 keep matching
# A30/E30 exactly.
$ws.Range("A30").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("E30").Copy()
$ws.Range("E31").PasteSpecial(-4122)

# --- Row 32: 2719. 统计整数数目 ------------------------------------------
$ws.Range("A30:E30").Copy($ws.Range("A32:E32"))

$ws.Range("A32").Value = "2719. 统计整数数目"
$ws.Range("B32").Value = "数学、动态规划"
$ws.Range("C32").Value = "⭐⭐"
$ws.Range("D32").Value = "2024.01.16"
$ws.Range("E32").Value = "力扣官方"

$ws.Hyperlinks.Add($ws.Range("A32"), "https://leetcode.cn/problems/count-of-integers/", "", "", "https://leetcode.cn/problems/count-of-integers/") | Out-Null
$ws.Range("A32").Value = "2719. 统计整数数目"
$ws.Hyperlinks.Add($ws.Range("E32"), "https://leetcode.cn/problems/count-of-integers/solutions/2613108/tong-ji-zheng-shu-shu-mu-by-leetcode-sol/?envType=daily-question&envId=2024-01-16") | Out-Null

$ws.Range("A30").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("E30").Copy()
$ws.Range("E32").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- View state: move the selection down to where the new rows are ------
$ws.Range("E34").Select()
